# Added Week 15 simulations
# - Adds a new RB player "J.Williams" (Week 15 simulation row, stats zeroed)
#   to the "RB" sheet as a new row right after the existing players.
# - Makes "RB" the active/selected sheet (previously "WR" was active),
#   with the selection left on cell K6 (first empty row's "extra" column).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("RB")

# New player row (row 5): Name + zeroed stat columns B:J
$ws.Range("A5").Value = "J.Williams"
$ws.Range("B5:J5").Value = 0

# Switch active sheet to RB and leave selection on K6, matching the saved view state
$ws.Activate()
$ws.Range("K6").Select()
